$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307, shifting existing rows 307-317 down to 308-318.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new data record.
$ws.Cells.Item(307, 1).Value = 4
$ws.Cells.Item(307, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(307, 3).Value = "Los Lagos"
$ws.Cells.Item(307, 4).Value = 44516
$ws.Cells.Item(307, 5).Value = 10
$ws.Cells.Item(307, 6).Value = "Fruta"
$ws.Cells.Item(307, 7).Value = 100103
$ws.Cells.Item(307, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(307, 9).Value = 100103006
$ws.Cells.Item(307, 10).Value = "Nectarín"
$ws.Cells.Item(307, 11).Value = "Early Glo"
$ws.Cells.Item(307, 12).Value = "Tercera"
$ws.Cells.Item(307, 13).Value = 400
$ws.Cells.Item(307, 14).Value = 13000
$ws.Cells.Item(307, 15).Value = 13500
$ws.Cells.Item(307, 16).Value = 13250
$ws.Cells.Item(307, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(307, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(307, 19).Value = 1893
$ws.Cells.Item(307, 20).Value = 7

# Give the new date cell the same date-time number format used by the
# other cells in column D (style index 2 in the original workbook).
$ws.Cells.Item(307, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
